# Update the "想去人数" (number of interested attendees) counts for three
# exhibition rows on both the "展览" sheet and the "全部类型" sheet.
#   F2: 286 -> 287
#   F3: 9   -> 10
#   F6: 212 -> 213

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 287
    $ws.Range("F3").Value = 10
    $ws.Range("F6").Value = 213
}
